$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow writing the updated values,
# then restore protection with the same password afterwards.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential footnote (row 38, col A)
# from 2021-05-25 to 2021-05-26.
$ws.Cells.Item(38, 1).Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-26 for illustrative purposes only and are subject to change."

# Refresh the Weight (col D) and Percent Change (col E) figures for each
# holding row (2-35) to the newly-supplied model values.
$ws.Range("D2").Value = 0.03543354750559932
$ws.Range("D3").Value = 0.02031474217684549
$ws.Range("E3").Value = 0.001953124999999778
$ws.Range("D4").Value = 0.01932993296074039
$ws.Range("E4").Value = 0.000400000000000178
$ws.Range("D5").Value = 0.03795788414553811
$ws.Range("E5").Value = -0.0006968641114982743
$ws.Range("D6").Value = 0.03420889451604715
$ws.Range("E6").Value = -0.000400000000000067
$ws.Range("D7").Value = 0.01976231304012538
$ws.Range("E7").Value = 0.00009652509652502417
$ws.Range("D8").Value = 0.03693924754675173
$ws.Range("E8").Value = 0.0004268943436498418
$ws.Range("D9").Value = 0.02048260738413613
$ws.Range("E9").Value = -0.0000894054537325939
$ws.Range("D10").Value = 0.02576730931911328
$ws.Range("E10").Value = -0.001974138781956403
$ws.Range("D11").Value = 0.02399073587528734
$ws.Range("E11").Value = 0.001855287569573205
$ws.Range("D12").Value = 0.05731528595838062
$ws.Range("E12").Value = 0.0009478672985780978
$ws.Range("D13").Value = 0.02490509257257499
$ws.Range("E13").Value = 0.0007352941176470562
$ws.Range("D14").Value = 0.02672902435215475
$ws.Range("E14").Value = -0.00561797752808979
$ws.Range("D15").Value = 0.0321780307172982
$ws.Range("E15").Value = -0.002124645892351174
$ws.Range("D16").Value = 0.01940628619593532
$ws.Range("E16").Value = -0.0104246122552758
$ws.Range("D17").Value = 0.03186895526290489
$ws.Range("E17").Value = 0.003671189146049469
$ws.Range("D18").Value = 0.04202327425664962
$ws.Range("E18").Value = -0.001149954001839704
$ws.Range("D19").Value = 0.1255390635054426
$ws.Range("E19").Value = 0
$ws.Range("D20").Value = 0.009067264605926252
$ws.Range("E20").Value = -0.003483870967741831
$ws.Range("D21").Value = 0.01521698104089654
$ws.Range("E21").Value = 0.001263601263601366
$ws.Range("D22").Value = 0.01767569764525809
$ws.Range("E22").Value = 0.002299412915851118
$ws.Range("D23").Value = 0.01561924711945848
$ws.Range("E23").Value = -0.007454739084132189
$ws.Range("D24").Value = 0.02176616580097392
$ws.Range("E24").Value = 0.003316249623153356
$ws.Range("D25").Value = 0.01258978881030875
$ws.Range("E25").Value = 0.009204114780725581
$ws.Range("D26").Value = 0.04250118141498161
$ws.Range("E26").Value = 0.00010771799429099
$ws.Range("D27").Value = 0.02386738038205104
$ws.Range("E27").Value = 0.00009803921568618534
$ws.Range("D28").Value = 0.04557174126579523
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0.05610218006036025
$ws.Range("E29").Value = 0.001066287542207167
$ws.Range("D30").Value = 0.01294484915785078
$ws.Range("E30").Value = 0.02413568166992808
$ws.Range("D31").Value = 0.0205550946327389
$ws.Range("E31").Value = -0.0003835826620637306
$ws.Range("D32").Value = 0.01355002866425585
$ws.Range("E32").Value = 0.000938526513373894
$ws.Range("D33").Value = 0.04167894710872527
$ws.Range("E33").Value = -0.0005162622612288059
$ws.Range("D34").Value = 0.01714122499889361
$ws.Range("E34").Value = -0.00290824487421848
$ws.Range("E35").Value = 0.0001497815463233909

$ws.Protect("D382")
